# Put logick to separate thread, using Task.Run().
# Fill in additional filled-form rows (4..15) on the ListSheet, mirroring the
# pattern already present in rows 2/3, and bump K3 (dn) from 3 to 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ListSheet")

# --- row 3: "dn" (K3) changes from 3 to 2 -----------------------------------
$ws.Range("K3").Value = 2

# --- clone row 3's formatting down to the new rows (4..15) -----------------
$ws.Range("A3:N3").Copy() | Out-Null
$ws.Range("A4:N15").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Columns that hold numeric-looking text (ids, dates, numeral degree) must be
# forced to Text so Excel doesn't reinterpret "14" as a number or
# "02.12.2017" as a date serial.
$ws.Range("A4:A15").NumberFormat = "@"
$ws.Range("D4:F15").NumberFormat = "@"
$ws.Range("L4:M15").NumberFormat = "@"

# --- data rows ---------------------------------------------------------------
$rows = @(
    @{ id = "14"; birth = "23.05.1997"; begin = "02.12.2017"; end = "03.04.2018"; fill = "04.04.2018"; hd = "182"; md = 47 },
    @{ id = "15"; birth = "23.05.1998"; begin = "02.12.2018"; end = "03.04.2019"; fill = "04.04.2019"; hd = "183"; md = 48 },
    @{ id = "16"; birth = "23.05.1999"; begin = "02.12.2019"; end = "03.04.2020"; fill = "04.04.2020"; hd = "184"; md = 49 },
    @{ id = "17"; birth = "23.05.2000"; begin = "02.12.2020"; end = "03.04.2021"; fill = "04.04.2021"; hd = "185"; md = 50 },
    @{ id = "18"; birth = "23.05.2001"; begin = "02.12.2021"; end = "03.04.2022"; fill = "04.04.2022"; hd = "186"; md = 51 },
    @{ id = "19"; birth = "23.05.2002"; begin = "02.12.2022"; end = "03.04.2023"; fill = "04.04.2023"; hd = "187"; md = 52 },
    @{ id = "20"; birth = "23.05.2003"; begin = "02.12.2023"; end = "03.04.2024"; fill = "04.04.2024"; hd = "188"; md = 53 },
    @{ id = "21"; birth = "23.05.2004"; begin = "02.12.2024"; end = "03.04.2025"; fill = "04.04.2025"; hd = "189"; md = 54 },
    @{ id = "22"; birth = "23.05.2005"; begin = "02.12.2025"; end = "03.04.2026"; fill = "04.04.2026"; hd = "190"; md = 55 },
    @{ id = "23"; birth = "23.05.2006"; begin = "02.12.2026"; end = "03.04.2027"; fill = "04.04.2027"; hd = "191"; md = 56 },
    @{ id = "24"; birth = "23.05.2007"; begin = "02.12.2027"; end = "03.04.2028"; fill = "04.04.2028"; hd = "192"; md = 57 },
    @{ id = "25"; birth = "23.05.2008"; begin = "02.12.2028"; end = "03.04.2029"; fill = "04.04.2029"; hd = "193"; md = 58 }
)

$r = 4
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row.id          # A - id
    $ws.Cells.Item($r, 2).Value = "DARYA"           # B - NAME
    $ws.Cells.Item($r, 3).Value = "BLABLABLA"       # C - LASTNAME
    $ws.Cells.Item($r, 4).Value = $row.birth        # D - birthdate
    $ws.Cells.Item($r, 5).Value = $row.begin        # E - begindate
    $ws.Cells.Item($r, 6).Value = $row.end          # F - enddate
    $ws.Cells.Item($r, 7).Value = "відмінно"        # G - ukrmark
    $ws.Cells.Item($r, 8).Value = "mit sehr gutem Erfolg"  # H - germark
    $ws.Cells.Item($r, 9).Value = "female"          # I - gender
    $ws.Cells.Item($r, 10).Value = "B"              # J - dl
    $ws.Cells.Item($r, 11).Value = 2                # K - dn
    $ws.Cells.Item($r, 12).Value = $row.fill        # L - filldate
    $ws.Cells.Item($r, 13).Value = $row.hd          # M - hd
    $ws.Cells.Item($r, 14).Value = $row.md          # N - md
    $r++
}

# --- final selection, as captured in the saved workbook ---------------------
$ws.Range("I18").Select() | Out-Null
